$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 168 (pushing old row 169 down to row 171)
$ws.Rows.Item(169).Insert()
$ws.Rows.Item(169).Insert()

# New row 169 = old row 167 values (before this edit)
$ws.Cells.Item(169, 1).Value = 9
$ws.Cells.Item(169, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(169, 3).Value = "Metropolitana"
$ws.Cells.Item(169, 4).Value = 44552
$ws.Cells.Item(169, 5).Value = 13
$ws.Cells.Item(169, 6).Value = "Fruta"
$ws.Cells.Item(169, 7).Value = 100101
$ws.Cells.Item(169, 8).Value = "Berries"
$ws.Cells.Item(169, 9).Value = 100101001
$ws.Cells.Item(169, 10).Value = "Arándano (blue)"
$ws.Cells.Item(169, 11).Value = "Sin especificar"
$ws.Cells.Item(169, 12).Value = "Primera"
$ws.Cells.Item(169, 13).Value = 310
$ws.Cells.Item(169, 14).Value = 3000
$ws.Cells.Item(169, 15).Value = 3000
$ws.Cells.Item(169, 16).Value = 3000
$ws.Cells.Item(169, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(169, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(169, 19).Value = 1500
$ws.Cells.Item(169, 20).Value = 2

# New row 170 = old row 168 values (before this edit)
$ws.Cells.Item(170, 1).Value = 9
$ws.Cells.Item(170, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(170, 3).Value = "Metropolitana"
$ws.Cells.Item(170, 4).Value = 44544
$ws.Cells.Item(170, 5).Value = 13
$ws.Cells.Item(170, 6).Value = "Fruta"
$ws.Cells.Item(170, 7).Value = 100101
$ws.Cells.Item(170, 8).Value = "Berries"
$ws.Cells.Item(170, 9).Value = 100101001
$ws.Cells.Item(170, 10).Value = "Arándano (blue)"
$ws.Cells.Item(170, 11).Value = "Sin especificar"
$ws.Cells.Item(170, 12).Value = "Primera"
$ws.Cells.Item(170, 13).Value = 650
$ws.Cells.Item(170, 14).Value = 3000
$ws.Cells.Item(170, 15).Value = 3000
$ws.Cells.Item(170, 16).Value = 3000
$ws.Cells.Item(170, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(170, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(170, 19).Value = 1500
$ws.Cells.Item(170, 20).Value = 2

# Apply the date number format (matching other D column cells) to the new D cells
$ws.Cells.Item(169, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(170, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Now update row 166 with the new values
$ws.Cells.Item(166, 4).Value = 44595
$ws.Cells.Item(166, 13).Value = 410

# Update row 167 with the new values
$ws.Cells.Item(167, 4).Value = 44595
$ws.Cells.Item(167, 13).Value = 330
$ws.Cells.Item(167, 14).Value = 3500
$ws.Cells.Item(167, 15).Value = 3500
$ws.Cells.Item(167, 16).Value = 3500
$ws.Cells.Item(167, 19).Value = 1750

# Update row 168 with the new values
$ws.Cells.Item(168, 4).Value = 44552
$ws.Cells.Item(168, 12).Value = "Especial"
$ws.Cells.Item(168, 13).Value = 280
$ws.Cells.Item(168, 14).Value = 4000
$ws.Cells.Item(168, 15).Value = 4000
$ws.Cells.Item(168, 16).Value = 4000
$ws.Cells.Item(168, 19).Value = 2000
